$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: rename the two shared strings in place on row 5 ---
# (write "Total Paths" before "Set Number" so the shared-string table keeps
#  the same allocation order as the target: index 9 = "Total Paths",
#  index 10 = "Set Number")
$ws.Range("A5").Value = "Total Paths"
$ws.Range("B5").Value = "Set Number"

# --- Phase 2: swap the contents of row 4 (old header) and row 5 (old set-label, just renamed) ---
$r4a = $ws.Range("A4").Value()
$r4b = $ws.Range("B4").Value()
$r4c = $ws.Range("C4").Value()
$r4d = $ws.Range("D4").Value()

$r5a = $ws.Range("A5").Value()
$r5b = $ws.Range("B5").Value()

# New row 4: "Set Number" / "Total Paths" (only two columns)
$ws.Range("A4").Value = $r5b
$ws.Range("B4").Value = $r5a
$ws.Range("C4").Clear()
$ws.Range("D4").Clear()

# New row 5: the header row that used to live in row 4
$ws.Range("A5").Value = $r4a
$ws.Range("B5").Value = $r4b
$ws.Range("C5").Value = $r4c
$ws.Range("D5").Value = $r4d

# --- Update the active selection to A5 ---
$ws.Range("A5").Select() | Out-Null
